{"js": "// Insert the literal marker \"##\" at a fixed set of points in the review\n// text (mostly right after a sentence-ending \". \"/\"... \" boundary, inside\n// the existing run so its formatting is preserved). Each entry below is an\n// (old, new) pair taken straight from the run's <w:t> text content; every\n// \"old\" string is unique in the document and fully contained in a single\n// run, so a body-wide search + in-place Replace keeps run/rPr formatting\n// untouched.\nconst replacements = [\n  [\". Ovaj film \", \". ## Ovaj film \"],\n  [\". Film je \", \". ## Film je \"],\n  [\n    \"(Hoking) kao i njegovo delo (Crne rupe). Intervjui sa \",\n    \"(Hoking) kao i njegovo delo (Crne rupe). ## Intervjui sa \",\n  ],\n  [\". Philip Glass \", \". ## Philip Glass \"],\n  [\n    \" film. Samo je jedan drugi \u010dovek mogao da komponuje takve proganjaju\u0107e \",\n    \" film. ## Samo je jedan drugi \u010dovek mogao da komponuje takve proganjaju\u0107e \",\n  ],\n  [\n    \" melodije (\u017dan Mi\u0161el Jarre). Sve u svemu bih visoko preporu\u010dio ovaj film \",\n    \" melodije (\u017dan Mi\u0161el Jarre). ## Sve u svemu bih visoko preporu\u010dio ovaj film \",\n  ],\n  [\" dugo... dugo vremena... \", \" dugo... dugo vremena... ## \"],\n  [\" posle 20 minuta... Keira \", \" posle 20 minuta... ## Keira \"],\n  [\" harizmu da popuni ulogu... \", \" harizmu da popuni ulogu... ## \"],\n  [\n    \": Da li je ikada imala \u010dasove glume? sude\u0107i po ivici ljubavi ona nikada nije bila u \",\n    \": Da li je ikada imala \u010dasove glume? ## sude\u0107i po ivici ljubavi ona nikada nije bila u \",\n  ],\n  [\" da ide u bliskoj budu\u0107nosti... \", \" da ide u bliskoj budu\u0107nosti... ## \"],\n  [\" u svojoj budu\u0107oj karijeri.. ako \", \" u svojoj budu\u0107oj karijeri.. ## ako \"],\n  [\". Sada uzimam ovo \", \". ## Sada uzimam ovo \"],\n  [\n    \" oko 4 meseca i anemija je nestala. Dobar proizvod. Lako se svari (za razliku od nekih drugih \",\n    \" oko 4 meseca i anemija je nestala. ## Dobar proizvod. ## Lako se svari (za razliku od nekih drugih \",\n  ],\n  [\n    \"Ovo je jedan od mojih omiljenih poslastica, i brzo se topi u ustima. Ovaj brend je dobar i isporu\u010duje \",\n    \"Ovo je jedan od mojih omiljenih poslastica, i brzo se topi u ustima. ## Ovaj brend je dobar i isporu\u010duje \",\n  ],\n  [\" dobro upakovan. Svako bi trebalo da \", \" dobro upakovan. ## Svako bi trebalo da \"],\n  [\" ovo jednom. \", \" ovo jednom. ## \"],\n  [\n    \"Ovo je fantasti\u010dna zagonetka/poklon za mlade i stare. To\",\n    \"Ovo je fantasti\u010dna zagonetka/poklon za mlade i stare. ## To\",\n  ],\n  [\" broj na\u010dina. \", \" broj na\u010dina. ## \"],\n  [\". U pore\u0111enju sa ve\u0107inom \", \". ## U pore\u0111enju sa ve\u0107inom \"],\n  [\" glavni lik je \u017eenka. \", \" glavni lik je \u017eenka. ## \"],\n  [\" je ista. Akcione scene nisu \", \" je ista. ## Akcione scene nisu \"],\n  [\". Specijalni efekti su tako \", \". ## Specijalni efekti su tako \"],\n  [\" naginjati na jednu stranu. Ote\u017eano \", \" naginjati na jednu stranu. ## Ote\u017eano \"],\n  [\". Imam malo iskustva\", \". ## Imam malo iskustva\"],\n  [\n    \"Ako pravilno obmotate zglobove, vide\u0107ete da su i suvi\u0161e uski i prekratki, suvi\u0161e kratki. Ne \",\n    \"Ako pravilno obmotate zglobove, vide\u0107ete da su i suvi\u0161e uski i prekratki, suvi\u0161e kratki. ## Ne \",\n  ],\n  [\n    \". Oni jednostavno ne\u0107e \u0161titi/podr\u017eati va\u0161e zglobove ili \",\n    \". ## Oni jednostavno ne\u0107e \u0161titi/podr\u017eati va\u0161e zglobove ili \",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  // Replacing the whole matched range (which equals one run's full text)\n  // in place keeps the existing run's formatting (rPr) intact.\n  found.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n\n// Special case: \"Bez obzira\" is its own (differently colored) run, with a\n// plain one-space run right before it. The diff turns that single space\n// into \" ## \" while leaving the \"Bez obzira\" run untouched. Searching for\n// \"Bez obzira\" and inserting \"## \" immediately before it achieves the same\n// resulting text, and the new text merges into the preceding plain run.\nconst bezResults = body.search(\"Bez obzira\", { matchCase: true });\nbezResults.load(\"items\");\nawait context.sync();\n\nif (bezResults.items.length === 0) {\n  throw new Error(\"No match found for: Bez obzira\");\n}\n\nbezResults.items[0].insertText(\"## \", Word.InsertLocation.before);\nawait context.sync();\n", "ps1": "# Insert the literal marker \"##\" at a fixed set of points in the review\n# text (mostly right after a sentence-ending \". \"/\"... \" boundary, inside\n# the existing run so its formatting is preserved). Each pair below is an\n# (old, new) run-text taken straight from the document; every \"old\" string\n# is unique in the document and fully contained in a single run, so a\n# whole-document Find/Replace keeps run/rPr formatting untouched.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @('. Ovaj film ', '. ## Ovaj film '),\n  @('. Film je ', '. ## Film je '),\n  @('(Hoking) kao i njegovo delo (Crne rupe). Intervjui sa ', '(Hoking) kao i njegovo delo (Crne rupe). ## Intervjui sa '),\n  @('. Philip Glass ', '. ## Philip Glass '),\n  @(' film. Samo je jedan drugi \u010dovek mogao da komponuje takve proganjaju\u0107e ', ' film. ## Samo je jedan drugi \u010dovek mogao da komponuje takve proganjaju\u0107e '),\n  @(' melodije (\u017dan Mi\u0161el Jarre). Sve u svemu bih visoko preporu\u010dio ovaj film ', ' melodije (\u017dan Mi\u0161el Jarre). ## Sve u svemu bih visoko preporu\u010dio ovaj film '),\n  @(' dugo... dugo vremena... ', ' dugo... dugo vremena... ## '),\n  @(' posle 20 minuta... Keira ', ' posle 20 minuta... ## Keira '),\n  @(' harizmu da popuni ulogu... ', ' harizmu da popuni ulogu... ## '),\n  @(': Da li je ikada imala \u010dasove glume? sude\u0107i po ivici ljubavi ona nikada nije bila u ', ': Da li je ikada imala \u010dasove glume? ## sude\u0107i po ivici ljubavi ona nikada nije bila u '),\n  @(' da ide u bliskoj budu\u0107nosti... ', ' da ide u bliskoj budu\u0107nosti... ## '),\n  @(' u svojoj budu\u0107oj karijeri.. ako ', ' u svojoj budu\u0107oj karijeri.. ## ako '),\n  @('. Sada uzimam ovo ', '. ## Sada uzimam ovo '),\n  @(' oko 4 meseca i anemija je nestala. Dobar proizvod. Lako se svari (za razliku od nekih drugih ', ' oko 4 meseca i anemija je nestala. ## Dobar proizvod. ## Lako se svari (za razliku od nekih drugih '),\n  @('Ovo je jedan od mojih omiljenih poslastica, i brzo se topi u ustima. Ovaj brend je dobar i isporu\u010duje ', 'Ovo je jedan od mojih omiljenih poslastica, i brzo se topi u ustima. ## Ovaj brend je dobar i isporu\u010duje '),\n  @(' dobro upakovan. Svako bi trebalo da ', ' dobro upakovan. ## Svako bi trebalo da '),\n  @(' ovo jednom. ', ' ovo jednom. ## '),\n  @('Ovo je fantasti\u010dna zagonetka/poklon za mlade i stare. To', 'Ovo je fantasti\u010dna zagonetka/poklon za mlade i stare. ## To'),\n  @(' broj na\u010dina. ', ' broj na\u010dina. ## '),\n  @('. U pore\u0111enju sa ve\u0107inom ', '. ## U pore\u0111enju sa ve\u0107inom '),\n  @(' glavni lik je \u017eenka. ', ' glavni lik je \u017eenka. ## '),\n  @(' je ista. Akcione scene nisu ', ' je ista. ## Akcione scene nisu '),\n  @('. Specijalni efekti su tako ', '. ## Specijalni efekti su tako '),\n  @(' naginjati na jednu stranu. Ote\u017eano ', ' naginjati na jednu stranu. ## Ote\u017eano '),\n  @('. Imam malo iskustva', '. ## Imam malo iskustva'),\n  @('Ako pravilno obmotate zglobove, vide\u0107ete da su i suvi\u0161e uski i prekratki, suvi\u0161e kratki. Ne ', 'Ako pravilno obmotate zglobove, vide\u0107ete da su i suvi\u0161e uski i prekratki, suvi\u0161e kratki. ## Ne '),\n  @('. Oni jednostavno ne\u0107e \u0161titi/podr\u017eati va\u0161e zglobove ili ', '. ## Oni jednostavno ne\u0107e \u0161titi/podr\u017eati va\u0161e zglobove ili ')\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $range = $d.Content\n  $found = $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n  if (-not $found) {\n    Write-Output \"MISSING: $old\"\n  }\n}\n\n# Special case: \"Bez obzira\" is its own (differently colored) run, with a\n# plain one-space run right before it. The diff turns that single space\n# into \" ## \" while leaving the \"Bez obzira\" run untouched. Finding\n# \"Bez obzira\" and inserting \"## \" immediately before its start achieves\n# the same resulting text, merging the new text into the preceding plain\n# run instead of touching \"Bez obzira\" itself.\n$bezRange = $d.Content\n$foundBez = $bezRange.Find.Execute(\"Bez obzira\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\nif ($foundBez) {\n  $insertPoint = $d.Range($bezRange.Start, $bezRange.Start)\n  $insertPoint.InsertBefore(\"## \")\n} else {\n  Write-Output \"MISSING: Bez obzira\"\n}\n"}
